# Daily attendance processing - 2026-01-18 14:59:09
# Swap the order of names in the "Recorded By" (column G) cells that
# currently read "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = $cell.Value()
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
